$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns based on the diff
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.571.14"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.596.94"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.24"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.56"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.66"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.052.16"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60.566.65"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.69"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000140"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.598.37"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "351.45"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.57"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.713.17"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0843"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +9.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.37"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +2.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.77"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.16"
$ws.Range("E35").Value = "  +4.30%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.916"
$ws.Range("E37").Value = "  +7.40%  "
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.37"
$ws.Range("E40").Value = "  +1.62%  "
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.53"
$ws.Range("E42").Value = "  -3.75%  "
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0558"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.53"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.01"
$ws.Range("E51").Value = "  +7.78%  "

# Rows 48/49: VeChain and RenderToken swap places in ranking order
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.86"
$ws.Range("E48").Value = "  +0.15%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0237"
$ws.Range("E49").Value = "  +1.33%  "
